# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that the file is now ready for handoff (new handoff xliff
# generated, content no longer flagged as duplicate, and a new error
# detail message about the handback file being out of date).

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06fda8377912ac043eda00a15b0edcff1537bd08/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bec0aa7d75a40c001e806a214e725504b157dd65/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row.
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady
$overview.Range("F3").Value = $statusReady
$overview.Range("G3").Value = "2016-09-07 00:52:43"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, I=Latest Target File,
#          J=Latest Handback File, K=Latest Handback DateTime,
#          L=Reference Tokens, M=To be localized, N=Dependency From,
#          O=Has metadata, P=Error Detail
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
# Leading apostrophe forces text so "False" isn't auto-converted to a Boolean;
# ClearFormats() afterwards drops the quote-prefix marker it leaves behind so
# the cell keeps its original (default) style.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").ClearFormats()
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 00:52:38"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.166666666666667

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row. Same column layout as zh-cn.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("F3").Value = "'False"
$dede.Range("F3").ClearFormats()
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 00:52:43"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.166666666666667
